$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (Coin name, Link, Volume) ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E49").Value = "48BOLOBOLO"

# --- Numeric-looking Price cells kept as text via scratch cell + PasteSpecial(values) ---
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "263.03"
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$scratch.Value = "22.89"
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$scratch.Value = "6.190"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$scratch.Value = "0.06243"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Value = "6.737"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Value = "3.449"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$scratch.Value = "1.346"
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Value = "0.7963"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Value = "0.1577"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Value = "0.08120"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$scratch.Value = "0.03415"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Value = "0.03085"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Value = "0.09334"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Value = "3.696"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Value = "0.001699"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$scratch.Value = "0.04811"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Value = "0.0006135"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$scratch.Value = "0.006235"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Value = "0.006181"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Value = "0.001093"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Value = "3.694"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Value = "2.216"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Value = "0.1274"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$scratch.Value = "0.0003199"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Value = "0.04638"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Value = "0.007072"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Value = "0.1122"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Value = "0.003146"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Value = "0.01017"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Value = "0.002968"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$scratch.Value = "0.00005883"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Value = "0.6995"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Value = "0.09030"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Value = "0.00002098"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
